# Changed header for child-parent tags report
#
# Removes the long list of childless "PUMP:*" tag paragraphs that used to
# follow the bold header paragraph, leaving only the title and the header
# line in place.

$d = $word.ActiveDocument

# The first two paragraphs are the document title ("Childless Report") and
# the bold header line ("These are the childless tags ..."); everything
# from the third paragraph through the end of the body is the list of
# childless tags that needs to be removed.
$count = $d.Paragraphs.Count
if ($count -gt 2) {
    $startPos = $d.Paragraphs(3).Range.Start
    $endPos = $d.Paragraphs($count).Range.End
    $listRange = $d.Range($startPos, $endPos)
    $listRange.Delete()
}
